# Update countries & provincias Spain
# Applies the refreshed COVID-19 "Pais" data dump:
#  - refreshed timestamp banner
#  - updated case statistics for many existing countries
#  - re-ranked country order for a few country groups (Irak/Catar,
#    Namibia/Benin/Sierra Leona/Ruanda, Lesoto/Taiwan/Vietnam,
#    Papua Nueva Guinea/Polinesia Francesa/... ) together with their
#    refreshed statistics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 26 de Julio de 2020 a las 15:48'
$ws.Cells.Item(4, 2).Value = 4317424
$ws.Cells.Item(4, 3).Value = 1715
$ws.Cells.Item(4, 5).Value = 2106325
$ws.Cells.Item(4, 7).Value = 9
$ws.Cells.Item(4, 8).Value = 149407
$ws.Cells.Item(6, 2).Value = 1412569
$ws.Cells.Item(6, 3).Value = 27075
$ws.Cells.Item(6, 4).Value = 902367
$ws.Cells.Item(6, 5).Value = 477843
$ws.Cells.Item(6, 7).Value = 263
$ws.Cells.Item(6, 8).Value = 32359
$ws.Cells.Item(16, 2).Value = 266941
$ws.Cells.Item(16, 3).Value = 1968
$ws.Cells.Item(16, 4).Value = 220323
$ws.Cells.Item(16, 5).Value = 43885
$ws.Cells.Item(16, 7).Value = 30
$ws.Cells.Item(16, 8).Value = 2733
$ws.Cells.Item(25, 1).Value = 'Irak'
$ws.Cells.Item(25, 2).Value = 110032
$ws.Cells.Item(25, 3).Value = 2459
$ws.Cells.Item(25, 4).Value = 75217
$ws.Cells.Item(25, 5).Value = 30453
$ws.Cells.Item(25, 7).Value = 78
$ws.Cells.Item(25, 8).Value = 4362
$ws.Cells.Item(26, 1).Value = 'Catar'
$ws.Cells.Item(26, 2).Value = 109305
$ws.Cells.Item(26, 3).Value = 269
$ws.Cells.Item(26, 4).Value = 106024
$ws.Cells.Item(26, 5).Value = 3116
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 165
$ws.Cells.Item(36, 2).Value = 67132
$ws.Cells.Item(36, 3).Value = 130
$ws.Cells.Item(36, 4).Value = 60425
$ws.Cells.Item(36, 5).Value = 6173
$ws.Cells.Item(36, 7).Value = 4
$ws.Cells.Item(36, 8).Value = 534
$ws.Cells.Item(39, 2).Value = 63773
$ws.Cells.Item(39, 3).Value = 464
$ws.Cells.Item(39, 4).Value = 54373
$ws.Cells.Item(39, 5).Value = 8967
$ws.Cells.Item(39, 7).Value = 4
$ws.Cells.Item(39, 8).Value = 433
$ws.Cells.Item(46, 2).Value = 50164
$ws.Cells.Item(46, 3).Value = 209
$ws.Cells.Item(46, 4).Value = 35217
$ws.Cells.Item(46, 5).Value = 13230
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = 1717
$ws.Cells.Item(57, 2).Value = 32437
$ws.Cells.Item(57, 3).Value = 586
$ws.Cells.Item(57, 4).Value = 28927
$ws.Cells.Item(57, 5).Value = 3349
$ws.Cells.Item(58, 2).Value = 30050
$ws.Cells.Item(58, 3).Value = 417
$ws.Cells.Item(58, 4).Value = 22684
$ws.Cells.Item(58, 5).Value = 6949
$ws.Cells.Item(58, 7).Value = 9
$ws.Cells.Item(58, 8).Value = 417
$ws.Cells.Item(62, 2).Value = 23730
$ws.Cells.Item(62, 3).Value = 467
$ws.Cells.Item(62, 5).Value = 9149
$ws.Cells.Item(62, 7).Value = 8
$ws.Cells.Item(62, 8).Value = 534
$ws.Cells.Item(82, 2).Value = 10086
$ws.Cells.Item(82, 3).Value = 152
$ws.Cells.Item(82, 4).Value = 5427
$ws.Cells.Item(82, 5).Value = 4199
$ws.Cells.Item(119, 2).Value = 2495
$ws.Cells.Item(119, 3).Value = 17
$ws.Cells.Item(119, 4).Value = 2349
$ws.Cells.Item(119, 5).Value = 59
$ws.Cells.Item(129, 2).Value = 1847
$ws.Cells.Item(129, 3).Value = 4
$ws.Cells.Item(129, 5).Value = 14
$ws.Cells.Item(130, 1).Value = 'Namibia'
$ws.Cells.Item(130, 2).Value = 1775
$ws.Cells.Item(130, 3).Value = 88
$ws.Cells.Item(130, 4).Value = 75
$ws.Cells.Item(130, 5).Value = 1692
$ws.Cells.Item(130, 8).Value = 8
$ws.Cells.Item(131, 1).Value = 'Benin'
$ws.Cells.Item(131, 2).Value = 1770
$ws.Cells.Item(131, 3).Value = 76
$ws.Cells.Item(131, 4).Value = 1036
$ws.Cells.Item(131, 5).Value = 699
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 35
$ws.Cells.Item(132, 1).Value = 'Sierra Leona'
$ws.Cells.Item(132, 2).Value = 1768
$ws.Cells.Item(132, 4).Value = 1297
$ws.Cells.Item(132, 5).Value = 405
$ws.Cells.Item(132, 8).Value = 66
$ws.Cells.Item(133, 1).Value = 'Ruanda'
$ws.Cells.Item(133, 2).Value = 1752
$ws.Cells.Item(133, 4).Value = 907
$ws.Cells.Item(133, 5).Value = 840
$ws.Cells.Item(133, 8).Value = 5
$ws.Cells.Item(159, 2).Value = 650
$ws.Cells.Item(159, 3).Value = 23
$ws.Cells.Item(159, 4).Value = 200
$ws.Cells.Item(159, 5).Value = 412
$ws.Cells.Item(159, 7).Value = 2
$ws.Cells.Item(159, 8).Value = 38
$ws.Cells.Item(161, 1).Value = 'Lesoto'
$ws.Cells.Item(161, 2).Value = 505
$ws.Cells.Item(161, 3).Value = 86
$ws.Cells.Item(161, 4).Value = 128
$ws.Cells.Item(161, 5).Value = 365
$ws.Cells.Item(161, 7).Value = 3
$ws.Cells.Item(161, 8).Value = 12
$ws.Cells.Item(162, 1).Value = 'Taiwan'
$ws.Cells.Item(162, 2).Value = 458
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 440
$ws.Cells.Item(162, 5).Value = 11
$ws.Cells.Item(162, 8).Value = 7
$ws.Cells.Item(163, 1).Value = 'Vietnam'
$ws.Cells.Item(163, 2).Value = 420
$ws.Cells.Item(163, 3).Value = 3
$ws.Cells.Item(163, 4).Value = 365
$ws.Cells.Item(163, 5).Value = 55
$ws.Cells.Item(163, 8).Value = 0
$ws.Cells.Item(192, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(192, 3).Value = 23
$ws.Cells.Item(192, 4).Value = 11
$ws.Cells.Item(192, 5).Value = 51
$ws.Cells.Item(193, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(193, 2).Value = 62
$ws.Cells.Item(193, 4).Value = 60
$ws.Cells.Item(193, 5).Value = 2
$ws.Cells.Item(194, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(194, 2).Value = 52
$ws.Cells.Item(194, 4).Value = 39
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 8).Value = 0
$ws.Cells.Item(195, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(195, 2).Value = 49
$ws.Cells.Item(195, 4).Value = 41
$ws.Cells.Item(195, 5).Value = 5
$ws.Cells.Item(195, 8).Value = 3
$ws.Cells.Item(196, 1).Value = 'Belice'
$ws.Cells.Item(196, 2).Value = 48
$ws.Cells.Item(196, 4).Value = 26
$ws.Cells.Item(196, 5).Value = 20
$ws.Cells.Item(196, 8).Value = 2
$ws.Cells.Item(197, 1).Value = 'Macao'
$ws.Cells.Item(197, 2).Value = 46
$ws.Cells.Item(197, 4).Value = 46
$ws.Cells.Item(197, 5).Value = 0
$ws.Cells.Item(197, 8).Value = 0
$ws.Cells.Item(198, 1).Value = 'Puerto Rico'
$ws.Cells.Item(198, 4).Value = 1
$ws.Cells.Item(198, 5).Value = 36
$ws.Cells.Item(198, 8).Value = 2
